$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1225
$ws.Range("I18").Value = 1378.3334
$ws.Range("J18").Value = 995
$ws.Range("K18").Value = 1378.3334
$ws.Range("L18").Value = 995
$ws.Range("M18").Value = -1094.3334
$ws.Range("N18").Value = -1563
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H98").Value = 1511.1852
$ws.Range("I98").Value = 1461.4286
$ws.Range("J98").Value = 1685.3334
$ws.Range("K98").Value = 1461.4286
$ws.Range("L98").Value = 1685.3334
$ws.Range("M98").Value = 36.57140000000004
$ws.Range("N98").Value = -4681.3334
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 738.125
$ws.Range("I107").Value = 738.125
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 738.125
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1181.875
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 1511.1852
$ws.Range("I122").Value = 1461.4286
$ws.Range("J122").Value = 1685.3334
$ws.Range("K122").Value = 4384.2858
$ws.Range("L122").Value = 5056.0002
$ws.Range("M122").Value = -1934.2858
$ws.Range("N122").Value = -9956.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2700
$ws.Range("I2").Value = 3050
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 3050
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -2937
$ws.Range("N2").Value = -2226
$ws.Range("H42").Value = 11500
$ws.Range("I42").Value = 3000
$ws.Range("K42").Value = 3000
$ws.Range("M42").Value = -2514
$ws.Range("H61").Value = 33400864
$ws.Range("I61").Value = 35750892
$ws.Range("J61").Value = 500507
$ws.Range("K61").Value = 35750892
$ws.Range("L61").Value = 500507
$ws.Range("M61").Value = -35750680
$ws.Range("N61").Value = -500931
$ws.Range("H74").Value = 7876121
$ws.Range("I74").Value = 13945218
$ws.Range("J74").Value = 72996.42999999999
$ws.Range("K74").Value = 13945218
$ws.Range("L74").Value = 72996.42999999999
$ws.Range("M74").Value = -13944344
$ws.Range("N74").Value = -74744.42999999999
$ws.Range("H77").Value = 7876121
$ws.Range("I77").Value = 13945218
$ws.Range("J77").Value = 72996.42999999999
$ws.Range("K77").Value = 69726090
$ws.Range("L77").Value = 364982.15
$ws.Range("M77").Value = -69721722
$ws.Range("N77").Value = -373718.15
$ws.Range("H110").Value = 5005005.5
$ws.Range("I110").Value = 5005005.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 5005005.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -5002960.5
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 2700
$ws.Range("I116").Value = 3050
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 3050
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = -756
$ws.Range("N116").Value = -6588
$ws.Range("H136").Value = 33400864
$ws.Range("I136").Value = 35750892
$ws.Range("J136").Value = 500507
$ws.Range("K136").Value = 107252676
$ws.Range("L136").Value = 1501521
$ws.Range("M136").Value = -107250126
$ws.Range("N136").Value = -1506621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2700
$ws.Range("I3").Value = 3050
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 3050
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -2936
$ws.Range("N3").Value = -2228
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H112").Value = 20000.334
$ws.Range("J112").Value = 20000.334
$ws.Range("L112").Value = 20000.334
$ws.Range("N112").Value = -22954.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 20835152
$ws.Range("J16").Value = 47620524
$ws.Range("L16").Value = 47620524
$ws.Range("N16").Value = -47621098
$ws.Range("H58").Value = 24880850
$ws.Range("I58").Value = 32906408
$ws.Range("J58").Value = 1619.8
$ws.Range("K58").Value = 32906408
$ws.Range("L58").Value = 1619.8
$ws.Range("M58").Value = -32906205
$ws.Range("N58").Value = -2025.8
$ws.Range("H81").Value = 43000
$ws.Range("J81").Value = 43000
$ws.Range("L81").Value = 43000
$ws.Range("N81").Value = -44996
$ws.Range("H84").Value = 43000
$ws.Range("J84").Value = 43000
$ws.Range("L84").Value = 129000
$ws.Range("N84").Value = -138984
$ws.Range("H110").Value = 44980
$ws.Range("J110").Value = 44980
$ws.Range("L110").Value = 44980
$ws.Range("N110").Value = -53160
$ws.Range("H113").Value = 20835152
$ws.Range("J113").Value = 47620524
$ws.Range("L113").Value = 47620524
$ws.Range("N113").Value = -47624864
$ws.Range("H122").Value = 1405.5483
$ws.Range("I122").Value = 1226.6207
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3679.8621
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -1229.8621
$ws.Range("N122").Value = -16900
$ws.Range("H136").Value = 24880850
$ws.Range("I136").Value = 32906408
$ws.Range("J136").Value = 1619.8
$ws.Range("K136").Value = 98719224
$ws.Range("L136").Value = 4859.4
$ws.Range("M136").Value = -98716674
$ws.Range("N136").Value = -9959.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 98085.71000000001
$ws.Range("I56").Value = 98085.71000000001
$ws.Range("K56").Value = 98085.71000000001
$ws.Range("M56").Value = -97555.71000000001
$ws.Range("H59").Value = 2000.5
$ws.Range("I59").Value = 2000.5
$ws.Range("K59").Value = 6001.5
$ws.Range("M59").Value = -5461.5
$ws.Range("H68").Value = 827.35
$ws.Range("J68").Value = 1051.5094
$ws.Range("L68").Value = 3154.5282
$ws.Range("N68").Value = -4776.5282
$ws.Range("H71").Value = 827.35
$ws.Range("J71").Value = 1051.5094
$ws.Range("L71").Value = 9463.584599999998
$ws.Range("N71").Value = -17575.5846
$ws.Range("H113").Value = 589.32355
$ws.Range("I113").Value = 512.375
$ws.Range("J113").Value = 613
$ws.Range("K113").Value = 1537.125
$ws.Range("L113").Value = 1839
$ws.Range("M113").Value = 632.875
$ws.Range("N113").Value = -6179
$ws.Range("H131").Value = 809.43475
$ws.Range("J131").Value = 951.0625
$ws.Range("L131").Value = 2853.1875
$ws.Range("N131").Value = -12933.1875
$ws.Range("H137").Value = 33684
$ws.Range("I137").Value = 1398
$ws.Range("J137").Value = 35477.668
$ws.Range("K137").Value = 4194
$ws.Range("L137").Value = 106433.004
$ws.Range("M137").Value = 906
$ws.Range("N137").Value = -116633.004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 54004.742
$ws.Range("I132").Value = 33608.72
$ws.Range("K132").Value = 100826.16
$ws.Range("M132").Value = -98296.16
$ws.Range("H136").Value = 30500
$ws.Range("J136").Value = 30500
$ws.Range("L136").Value = 91500
$ws.Range("N136").Value = -96600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2128.4666
$ws.Range("I7").Value = 2032.9231
$ws.Range("K7").Value = 2032.9231
$ws.Range("M7").Value = -1920.9231
$ws.Range("H126").Value = 2128.4666
$ws.Range("I126").Value = 2032.9231
$ws.Range("K126").Value = 6098.7693
$ws.Range("M126").Value = -3628.7693
$ws.Range("H127").Value = 49966.668
$ws.Range("J127").Value = 49966.668
$ws.Range("L127").Value = 49966.668
$ws.Range("N127").Value = -59886.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 4166.6665
$ws.Range("H57").Value = 52000
$ws.Range("I57").Value = 52000
$ws.Range("K57").Value = 52000
$ws.Range("M57").Value = -51246
$ws.Range("H113").Value = 2139.3076
$ws.Range("I113").Value = 727.875
$ws.Range("J113").Value = 4397.6
$ws.Range("K113").Value = 2183.625
$ws.Range("L113").Value = 13192.8
$ws.Range("M113").Value = -13.625
$ws.Range("N113").Value = -17532.8
$ws.Range("H136").Value = 208331.4
$ws.Range("I136").Value = 209200.8
$ws.Range("J136").Value = 207462
$ws.Range("K136").Value = 627602.3999999999
$ws.Range("L136").Value = 622386
$ws.Range("M136").Value = -625052.3999999999
$ws.Range("N136").Value = -627486
